# Annotated w/ BLAST info
# Adds BLAST-match annotation columns (F:J) to the top rows of the
# WinterContentABX_TaxTable sheet: species match, max score, query
# cover %, e-value and % identity for the top 6 OTUs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values -----------------------------------------------------
# NOTE: write order matters for shared-string de-dup indices, so this
# mirrors the order the values were actually typed in: species name
# down column F first, then the header labels (scattered), then the
# remaining species names, and finally the "BLAST Match" header.

$ws.Range("F2").Value = "Morganella morganii"

$ws.Range("G1").Value = "Max Score"
$ws.Range("I1").Value = "E-value"
$ws.Range("H1").Value = "Query Cover (%)"
$ws.Range("J1").Value = "% Identity"

$ws.Range("F3").Value = "Uncultured bacterium"
$ws.Range("F4").Value = "Latilactobacillus sakei / curvatus"
$ws.Range("F5").Value = "Alistipes shahii"
$ws.Range("F6").Value = "Uncultured Bacilli / Clostridium / bacterium"

$ws.Range("F1").Value = "BLAST Match"

$ws.Range("F7").Value = "Uncultured bacterium"

# --- Numeric columns (Max Score / Query Cover / E-value / % Identity) -
$ws.Range("G2").Value = 568
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = [double]"2E-127"
$ws.Range("J2").Value = 100

$ws.Range("G3").Value = 446
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = [double]"8E-121"
$ws.Range("J3").Value = 98.42

$ws.Range("G4").Value = 468
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = [double]"2E-127"
$ws.Range("J4").Value = 100

$ws.Range("G5").Value = 418
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = [double]"2E-112"
$ws.Range("J5").Value = 96.44

$ws.Range("G6").Value = 405
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = [double]"1E-108"
$ws.Range("J6").Value = 95.63

$ws.Range("G7").Value = 416
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = [double]"6E-112"
$ws.Range("J7").Value = 96.43

# --- Formatting --------------------------------------------------------
# G/H/J share the existing highlight fill used by C:E on these rows.
$highlightRange = $ws.Range("G2:H7")
$highlightRange.Interior.Color = 65535
$ws.Range("J2:J7").Interior.Color = 65535

# E-value column gets scientific-notation formatting.
$ws.Range("I2:I7").NumberFormat = "0.00E+00"

# --- Column widths for the new columns ---------------------------------
$ws.Columns.Item(6).ColumnWidth = 19.5703125
$ws.Columns.Item(7).ColumnWidth = 10.140625
$ws.Columns.Item(8).ColumnWidth = 15.5703125
$ws.Columns.Item(9).ColumnWidth = 9.28515625
$ws.Columns.Item(10).ColumnWidth = 15.42578125

# --- Selection moves to the first newly-entered cell --------------------
[void]$ws.Range("F2").Select()
